$d = $word.ActiveDocument

# --- 1) Merge the two "Addendum change to Roll(" / ") in Die Class (Die.cs)"
#        runs (with the _GoBack bookmark collapsed between them) into a
#        single run, dropping the bookmark from this spot. A find/replace
#        across the whole phrase (which spans the old run + bookmark +
#        run boundary) collapses it into one run and removes the
#        now-interior bookmark.
$rHeading = $d.Content
$rHeading.Find.Execute("Addendum change to Roll() in Die Class (Die.cs)", `
    $false, $false, $false, $false, $false, $true, 0, $false, `
    "Addendum change to Roll() in Die Class (Die.cs)", 2) | Out-Null

# --- 2) Re-insert the _GoBack bookmark, collapsed, right before the first
#        run of the "Changes to Roll() in Die.cs implemented" paragraph.
$rTarget = $d.Content
$rTarget.Find.Execute("Changes to Roll() in Die.cs implemented", `
    $false, $false, $false, $false, $false, $true, 0, $false, "", 0) | Out-Null
$rBookmark = $d.Range($rTarget.Start, $rTarget.Start)
$d.Bookmarks.Add("_GoBack", $rBookmark) | Out-Null

# --- 3) Flip the two "spec change" answers from Y to N.
$rClass = $d.Content
$rClass.Find.Execute("Did you make any changes to the specification of any class", `
    $false, $false, $false, $false, $false, $true, 0, $false, "", 0) | Out-Null
$rClassAnswer = $d.Range($rClass.End, $rClass.End + 40)
$rClassAnswer.Find.Execute("Y", $true, $false, $false, $false, $false, `
    $true, 0, $false, "N", 1) | Out-Null

$rMethod = $d.Content
$rMethod.Find.Execute("Did you make any changes to the specification of any supplied method", `
    $false, $false, $false, $false, $false, $true, 0, $false, "", 0) | Out-Null
$rMethodAnswer = $d.Range($rMethod.End, $rMethod.End + 40)
$rMethodAnswer.Find.Execute("Y", $true, $false, $false, $false, $false, `
    $true, 0, $false, "N", 1) | Out-Null
